$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.203.09"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "2.143.01"
$ws.Range("E3").Value = "  -3.37%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.45"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.599"
$ws.Range("E6").Value = "  -4.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.97"
$ws.Range("E7").Value = "  -5.71%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -7.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.20"
$ws.Range("E10").Value = "  -10.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0894"
$ws.Range("E11").Value = "  -6.86%  "
$ws.Range("E12").Value = "  -7.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0995"
$ws.Range("E13").Value = "  -4.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.58"
$ws.Range("E14").Value = "  -6.97%  "
$ws.Range("D15").Value = "2.460.96"
$ws.Range("E15").Value = "  -3.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.38"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "2.163.79"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.774"
$ws.Range("E18").Value = "  -7.72%  "
$ws.Range("D19").Value = "41.029.71"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "0.0₃0995"
$ws.Range("E20").Value = "  -7.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.10"
$ws.Range("E21").Value = "  -5.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.70"
$ws.Range("E22").Value = "  -8.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.45"
$ws.Range("E23").Value = "  -12.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "223.67"
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.87"
$ws.Range("E26").Value = "  -10.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.52"
$ws.Range("E27").Value = "  -11.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.36"
$ws.Range("E28").Value = "  -8.48%  "
$ws.Range("E29").Value = "  -6.11%  "
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.48"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.54"
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.64"
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0748"
$ws.Range("E34").Value = "  -6.04%  "
$ws.Range("E35").Value = "  -11.73%  "
$ws.Range("E36").Value = "  -5.16%  "
$ws.Range("E37").Value = "  -9.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.06"
$ws.Range("E38").Value = "  -5.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.02"
$ws.Range("E40").Value = "  -4.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.58"
$ws.Range("E41").Value = "  -16.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.24"
$ws.Range("E42").Value = "  -7.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "57.08"
$ws.Range("E43").Value = "  -12.98%  "
$ws.Range("E44").Value = "  -7.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.12"
$ws.Range("E45").Value = "  -7.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0946"
$ws.Range("E46").Value = "  -6.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "96.27"
$ws.Range("E47").Value = "  -8.92%  "
$ws.Range("E48").Value = "  -4.99%  "
$ws.Range("E49").Value = "  -6.25%  "
$ws.Range("E50").Value = "  -12.72%  "
